$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need to be forced
# to remain stored as text (matching the original inlineStr cells), since
# Excel would otherwise auto-convert them to numeric values on assignment.

$ws.Range('D2').Value = '30.062.97'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '2.104.99'
$ws.Range('E3').Value = '  -0.56%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.70%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '347.69'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +3.10%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.69%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5168'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -1.51%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.4444'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.66%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '52.30'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -4.18%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.08976'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.88%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '1.173'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.15%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '25.46'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +3.63%  '
$ws.Range('D13').Value = '2.110.16'
$ws.Range('E13').Value = '  -0.44%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '8.254'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +1.53%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '6.728'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -2.07%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '99.30'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +2.12%  '
$ws.Range('E17').Value = '  -2.30%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.62%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '20.86'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +7.34%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.06678'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  -0.59%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.241'
$c.Style = 'Normal'
$ws.Range('D23').Value = '30.154.77'
$ws.Range('E23').Value = '  -1.64%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '12.73'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -1.12%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.347'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').Value = '2.358.43'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '21.98'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.99%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.543'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.30%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '162.47'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.97%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '133.59'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.48%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.177'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -3.31%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.1066'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.69%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.639'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.85%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '6.242'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.16%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.960'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.11%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '5.921'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('E37').Value = '  -3.46%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.02577'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -2.01%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.06802'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.77%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.2302'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -1.07%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '12.57'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.42%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.6818'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.08%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.292'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.66%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '14.23'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -4.31%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.6388'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -1.16%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.294'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.26%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.00000000365'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.23%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.639'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('E49').Value = '  -2.80%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '82.51'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -1.26%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.07228'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.35%  '
